$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content fixes (DSM exporter "parent.child" label fix) ---
$ws.Range("C4").Value = "Component 1.1.Component 1.1.1"
$ws.Range("D4").Value = "Component 1.Component 1.2"
$ws.Range("E4").Value = "testReqTrace.Component 2"
$ws.Range("F4").Value = "testReqTrace.Component 3"
$ws.Range("G4").Value = "testReqTrace.Component 4"
$ws.Range("H4").Value = "testReqTrace.Component 6"

$ws.Range("A6").Value = ".testReqTraceUseCases"
$ws.Range("A7").Value = "testReqTraceUseCases.Requirement 1"
$ws.Range("A8").Value = "testReqTraceUseCases.Requirement 2"
$ws.Range("A9").Value = "testReqTraceUseCases.Requirement 3"
$ws.Range("A10").Value = "testReqTraceUseCases.Requirement 4"
$ws.Range("A11").Value = "testReqTraceUseCases.Requirement 5"
$ws.Range("A12").Value = "testReqTraceUseCases.Requirement 6"
$ws.Range("A13").Value = "testReqTraceUseCases.Requirement 7"

# --- Remove the leftover red/green highlight fills on the rotated
#     component headers (row 4) by matching the already-unfilled
#     neighbour cells G4/H4 -- keeps the 90 degree text rotation. ---
$ws.Range("G4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("F4").PasteSpecial(-4122)

# --- Remove the stray blue highlight fill from the requirement rows
#     that shouldn't have had one (match plain rows A7/A9/A10). ---
$ws.Range("A9").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)

# --- Slightly narrower first column ---
$ws.Columns.Item(1).ColumnWidth = 38.377604166666664
